$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Quantity changes (column G)
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G8").Value = 2

# Total Price column (H) = Unit Price (F) * Quantity (G) for every BOM row
$ws.Range("H2:H18").Formula = "=F2*G2"

# Update the active selection to match the post-edit state
$ws.Range("E10").Select()
